# EPBDS-10709 - update EPBDS-8671_CollectParametersValidation.xlsx test data
#
# Updates the expected-result columns (E/F) of the first "Collect" test
# table (rows 4-7) and appends five new test-case rows (8-12) with a new
# "testX" test case, in line with the removal of
# SyntaxNodeExceptionCollector.addSyntaxNodeException() usage.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A helper cell that carries the sheet's plain/default style (no number
# format, no quote-prefix) so we can restore it onto cells after we force
# a numeric-looking literal to be stored as text.
$plainStyleCell = $ws.Cells.Item(3, 2)

function Set-TextValue($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.Value = "'" + $text
    $cell.Style = $plainStyleCell.Style
}

# ---------------------------------------------------------------------
# Update existing rows 4-7 (expected-result columns E/F)
# ---------------------------------------------------------------------
Set-TextValue 4 6 "27"    # F4
Set-TextValue 5 5 "2.6"   # E5
Set-TextValue 6 5 "25"    # E6
Set-TextValue 6 6 "24"    # F6
Set-TextValue 7 5 "23"    # E7
Set-TextValue 7 6 "2.2"   # F7

# ---------------------------------------------------------------------
# New rows 8-12: additional "testX" test case data
# ---------------------------------------------------------------------
Set-TextValue 8 2 "2"       # B8
$ws.Cells.Item(8, 3).Value = "testX"   # C8
Set-TextValue 8 4 "3"       # D8
Set-TextValue 8 5 "4"       # E8

Set-TextValue 9 2 "5"       # B9
$ws.Cells.Item(9, 3).Value = "testX"   # C9
Set-TextValue 9 4 "6"       # D9
$ws.Cells.Item(9, 6).Value = 7         # F9 (numeric)

Set-TextValue 10 2 "8"      # B10
$ws.Cells.Item(10, 3).Value = "testX"  # C10
Set-TextValue 10 4 "9"      # D10

Set-TextValue 11 2 "10"     # B11
$ws.Cells.Item(11, 3).Value = "testX"  # C11
Set-TextValue 11 5 "11"     # E11

Set-TextValue 12 2 "12"     # B12
$ws.Cells.Item(12, 3).Value = "testX"  # C12
